$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 16; this shifts the existing rows 16..27 down
# to 17..28 (carrying their values/formatting along), making room at row 16
# for a new weekly price record.
$ws.Rows(16).Insert()

# Populate the newly inserted row 16 with the new weekly record.
$ws.Cells.Item(16, 1).Value = 8
$ws.Cells.Item(16, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(16, 3).Value = "Coquimbo"
$ws.Cells.Item(16, 4).Value = 44651
$ws.Cells.Item(16, 4).NumberFormat = $ws.Cells.Item(17, 4).NumberFormat
$ws.Cells.Item(16, 5).Value = 4
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100101
$ws.Cells.Item(16, 8).Value = "Berries"
$ws.Cells.Item(16, 9).Value = 100101001
$ws.Cells.Item(16, 10).Value = "Arándano (blue)"
$ws.Cells.Item(16, 11).Value = "Sin especificar"
$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value = 400
$ws.Cells.Item(16, 14).Value = 6000
$ws.Cells.Item(16, 15).Value = 6500
$ws.Cells.Item(16, 16).Value = 6250
$ws.Cells.Item(16, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(16, 18).Value = "Provincia de Linares"
$ws.Cells.Item(16, 19).Value = 3125
$ws.Cells.Item(16, 20).Value = 2
